$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename header columns (row 1) ---
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# --- Title-case the "de"/"del"/"el" connector words in a handful of place names ---
$ws.Range("A17").Value  = "Ciudad De México"
$ws.Range("A28").Value  = "Estado De México"
$ws.Range("B29").Value  = "Ecatepec De Morelos"
$ws.Range("B32").Value  = "Tlalnepantla De Baz"
$ws.Range("B34").Value  = "Apaseo El Grande"
$ws.Range("B41").Value  = "Acapulco De Juárez"
$ws.Range("B47").Value  = "Tula De Allende"
$ws.Range("B51").Value  = "Lagos De Moreno"
$ws.Range("B53").Value  = "Tepatitlán De Morelos"
$ws.Range("B75").Value  = "Nejapa De Madero"
$ws.Range("B88").Value  = "San Juan Del Río"
$ws.Range("B93").Value  = "Villa De Ramos"
$ws.Range("B108").Value = "Nanacamilpa De Mariano Arista"

# --- Remove the trailing metadata/footer rows (120-124), shifting the dimension to A1:D118 ---
$ws.Rows("120:124").Delete()
